$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.343.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.76%  "

$ws.Range("D3").Value = "'3.616.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.68%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'627.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.26%  "

$ws.Range("D6").Value = "'159.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.49%  "

$ws.Range("D7").Value = "'3.616.28"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.65%  "

$ws.Range("D9").Value = "'0.492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.40%  "

$ws.Range("E10").Value = "  +5.60%  "

$ws.Range("D11").Value = "'7.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.27%  "

$ws.Range("D12").Value = "'0.439"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.57%  "

$ws.Range("E13").Value = "  +3.06%  "

$ws.Range("D14").Value = "'33.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.74%  "

$ws.Range("D15").Value = "'4.228.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.79%  "

$ws.Range("D16").Value = "'3.615.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.29%  "

$ws.Range("D17").Value = "'69.214.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.57%  "

$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("E19").Value = "  +6.34%  "

$ws.Range("D20").Value = "'15.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.34%  "

$ws.Range("D21").Value = "'10.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.98%  "

$ws.Range("D22").Value = "'459.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.07%  "

$ws.Range("E23").Value = "  +2.57%  "

$ws.Range("D24").Value = "'78.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.42%  "

$ws.Range("E25").Value = "  +13.57%  "

$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "'3.760.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.73%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'10.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.56%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").Value = "'9.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +13.13%  "

$ws.Range("D30").Value = "'2.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.61%  "

$ws.Range("D31").Value = "'1.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.60%  "

$ws.Range("E32").Value = "  +12.07%  "

$ws.Range("D33").Value = "'6.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.07%  "

$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("E35").Value = "  +6.81%  "

$ws.Range("E36").Value = "  +3.51%  "

$ws.Range("D37").Value = "'3.610.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.70%  "

$ws.Range("D38").Value = "'8.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.03%  "

$ws.Range("D39").Value = "'2.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.55%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").Value = "'0.0925"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.32%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").Value = "'175.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.20%  "

$ws.Range("D44").Value = "'5.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.92%  "

$ws.Range("D45").Value = "'31.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +17.59%  "

$ws.Range("E46").Value = "  +2.54%  "

$ws.Range("E47").Value = "  +13.14%  "

$ws.Range("D48").Value = "'2.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.08%  "

$ws.Range("D49").Value = "'46.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.21%  "

$ws.Range("D50").Value = "'7.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.76%  "

$ws.Range("E51").Value = "  +7.86%  "
